$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the discrete-time execution values (column I) per the corrected
# Matlab/SciLab simulation data.
$ws.Range("I2").Value = 1.8541993999999999
$ws.Range("I3").Value = 2.0832563999999998
$ws.Range("I4").Value = 2.0191528999999999
$ws.Range("I7").Value = 4.0159577000000004
$ws.Range("I10").Value = 2.9229202000000001
$ws.Range("I13").Value = 3.0300665000000002
$ws.Range("I14").Value = 2.9970636000000002
$ws.Range("I16").Value = 3.1124489
$ws.Range("I17").Value = 2.7545639

# Row 4 now references "LoopAlg pero bien" instead of the (removed) "mala"
# comment; rows 7 and 17 no longer carry the erroneous "mala" annotation.
$ws.Range("J4").Value = "LoopAlg pero bien"
$ws.Range("J7").ClearContents()
$ws.Range("J17").ClearContents()

# Move the active selection to J5, matching the saved view state.
$ws.Range("J5").Select()
